# Update Combo.xlsx: bảng combo và hàng tặng hàng
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data row (row 2): clear old per-cell formatting first so the cells
#     fall back to the default style (no border / no alignment / no text
#     number-format), matching the refreshed data. ---
$ws.Range("A2:E2").ClearFormats()

# --- Header row (row 1): clear old formatting too, we rebuild it below. ---
$ws.Range("A1:E1").ClearFormats()

# --- Header row values ---
$ws.Range("A1").Value = "Nhóm khách hàng"
$ws.Range("B1").Value = "Mã sản phẩm"
$ws.Range("C1").Value = "Tên sản phẩm"
$ws.Range("D1").Value = "Mã Barcode"
$ws.Range("E1").Value = "Trạng thái"

# --- Data row values ---
$ws.Range("A2").Value = "Emart"
$ws.Range("B2").Value = 50011840
$ws.Range("C2").Value = "Bút vẽ lên vải FM-C002 túi 12 màu"
$ws.Range("D2").Value = 8935001868620
$ws.Range("E2").Value = "x"

# --- Header formatting: bold font, black border, light-blue fill, centered ---
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 0
$headerRange.Interior.Color = 14599344
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.Borders.Color = 0
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 15.1666666
$ws.Columns.Item(2).ColumnWidth = 11.1666666
$ws.Columns.Item(3).ColumnWidth = 59.1666666
$ws.Columns.Item(4).ColumnWidth = 13.1666666
$ws.Columns.Item(5).ColumnWidth = 10.1666666

# --- Selection ---
$ws.Range("H9").Select() | Out-Null
